$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "Andre Lucca-Circuitos Elétricos"
$ws.Range("F4").Value = "Andre Lucca-Circuitos Elétricos"
$ws.Range("D6").Value = "José Ferreira-Tecnologia dos Materiais"
$ws.Range("E6").Value = "-"
$ws.Range("E7").Value = "-"
